$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format to prevent Excel from auto-converting
# numeric-looking strings (e.g. "235.15") into actual numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value2 = '92.548.74'
$ws.Range("E2").Value2 = '  +0.74%  '
$ws.Range("D3").Value2 = '3.110.03'
$ws.Range("E3").Value2 = '  -0.65%  '
$ws.Range("E4").Value2 = '  +0.01%  '
$ws.Range("D5").Value2 = '235.15'
$ws.Range("E5").Value2 = '  -2.85%  '
$ws.Range("D6").Value2 = '612.62'
$ws.Range("E6").Value2 = '  -0.97%  '
$ws.Range("E7").Value2 = '  -1.53%  '
$ws.Range("D8").Value2 = '0.388'
$ws.Range("E8").Value2 = '  -0.04%  '
$ws.Range("E9").Value2 = '  -0.04%  '
$ws.Range("B10").Value2 = 'Cardano'
$ws.Range("C10").Value2 = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value2 = '0.793'
$ws.Range("E10").Value2 = '  +5.14%  '
$ws.Range("B11").Value2 = 'LidoStakedEther'
$ws.Range("C11").Value2 = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D11").Value2 = '3.107.18'
$ws.Range("E11").Value2 = '  -0.67%  '
$ws.Range("E12").Value2 = '  -4.03%  '
$ws.Range("E13").Value2 = '  -3.63%  '
$ws.Range("D14").Value2 = '92.277.83'
$ws.Range("E14").Value2 = '  +0.81%  '
$ws.Range("D15").Value2 = '33.95'
$ws.Range("E15").Value2 = '  -4.03%  '
$ws.Range("D16").Value2 = '5.41'
$ws.Range("E16").Value2 = '  -3.60%  '
$ws.Range("D17").Value2 = '3.693.42'
$ws.Range("E17").Value2 = '  -0.48%  '
$ws.Range("D18").Value2 = '3.096.05'
$ws.Range("E18").Value2 = '  -1.20%  '
$ws.Range("D19").Value2 = '3.81'
$ws.Range("E19").Value2 = '  +0.23%  '
$ws.Range("D20").Value2 = '14.54'
$ws.Range("E20").Value2 = '  -3.18%  '
$ws.Range("D21").Value2 = '5.85'
$ws.Range("E21").Value2 = '  -1.32%  '
$ws.Range("E22").Value2 = '  +0.20%  '
$ws.Range("B23").Value2 = 'BitcoinCash'
$ws.Range("C23").Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value2 = '439.68'
$ws.Range("E23").Value2 = '  -4.08%  '
$ws.Range("B24").Value2 = 'Uniswap'
$ws.Range("C24").Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").Value2 = '9.13'
$ws.Range("E24").Value2 = '  -1.37%  '
$ws.Range("B25").Value2 = 'LEO'
$ws.Range("C25").Value2 = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").Value2 = '8.13'
$ws.Range("E25").Value2 = '  +4.52%  '
$ws.Range("B26").Value2 = 'NEARProtocol'
$ws.Range("C26").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value2 = '5.58'
$ws.Range("E26").Value2 = '  -6.34%  '
$ws.Range("B27").Value2 = 'Litecoin'
$ws.Range("C27").Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value2 = '85.54'
$ws.Range("E27").Value2 = '  -4.44%  '
$ws.Range("B28").Value2 = 'Aptos'
$ws.Range("C28").Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value2 = '11.52'
$ws.Range("E28").Value2 = '  -2.03%  '
$ws.Range("B29").Value2 = 'WrappedeETH'
$ws.Range("C29").Value2 = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value2 = '3.272.83'
$ws.Range("E29").Value2 = '  -0.58%  '
$ws.Range("B30").Value2 = 'Dai'
$ws.Range("C30").Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value2 = '1.00'
$ws.Range("E30").Value2 = '  +0.10%  '
$ws.Range("B31").Value2 = 'Cronos'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value2 = '0.180'
$ws.Range("E31").Value2 = '  +7.26%  '
$ws.Range("B32").Value2 = 'Stellar'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value2 = '0.235'
$ws.Range("E32").Value2 = '  +4.13%  '
$ws.Range("B33").Value2 = 'Hedera'
$ws.Range("C33").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value2 = '0.123'
$ws.Range("E33").Value2 = '  -15.20%  '
$ws.Range("D34").Value2 = '1.04'
$ws.Range("E34").Value2 = '  +45.95%  '
$ws.Range("B35").Value2 = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value2 = '9.18'
$ws.Range("E35").Value2 = '  -2.93%  '
$ws.Range("B36").Value2 = 'RenderToken'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value2 = '8.14'
$ws.Range("E36").Value2 = '  +8.40%  '
$ws.Range("B37").Value2 = 'Kaspa'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value2 = '0.162'
$ws.Range("E37").Value2 = '  -8.14%  '
$ws.Range("B38").Value2 = 'EthereumClassic'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value2 = '25.69'
$ws.Range("E38").Value2 = '  -3.28%  '
$ws.Range("B39").Value2 = 'MantraDAO'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D39").Value2 = '3.96'
$ws.Range("E39").Value2 = '  +2.14%  '
$ws.Range("B40").Value2 = 'PancakeSwap'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D40").Value2 = '1.89'
$ws.Range("E40").Value2 = '  -3.34%  '
$ws.Range("B41").Value2 = 'WhiteBITCoin'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value2 = '23.87'
$ws.Range("E41").Value2 = '  +7.65%  '
$ws.Range("B42").Value2 = 'Fetch.AI'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value2 = '1.28'
$ws.Range("E42").Value2 = '  -3.27%  '
$ws.Range("B43").Value2 = 'Bittensor'
$ws.Range("C43").Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value2 = '464.44'
$ws.Range("E43").Value2 = '  -5.63%  '
$ws.Range("B44").Value2 = 'PolygonEcosystemToken'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").Value2 = '0.429'
$ws.Range("E44").Value2 = '  -2.63%  '
$ws.Range("B45").Value2 = 'dogwifhat'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value2 = '3.31'
$ws.Range("E45").Value2 = '  -2.53%  '
$ws.Range("B46").Value2 = 'USDe'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value2 = '1.00'
$ws.Range("E46").Value2 = '  +0.01%  '
$ws.Range("B47").Value2 = 'Monero'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").Value2 = '159.93'
$ws.Range("E47").Value2 = '  +2.13%  '
$ws.Range("B48").Value2 = 'ARBITRUM'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value2 = '0.682'
$ws.Range("E48").Value2 = '  -4.22%  '
$ws.Range("B49").Value2 = 'Stacks'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value2 = '1.83'
$ws.Range("E49").Value2 = '  -4.98%  '
$ws.Range("B50").Value2 = 'ImmutableX'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").Value2 = '1.33'
$ws.Range("E50").Value2 = '  -2.29%  '
$ws.Range("B51").Value2 = 'VeChain'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value2 = '0.0325'
$ws.Range("E51").Value2 = '  -0.77%  '

# Restore default (General) styling on column D now that text values are set,
# so the saved cells do not carry a stray explicit style index.
$dRange.ClearFormats()
